$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task description text (row 5 / col B): struct & menu description ---
$b5 = "*a Khai báo struct cho sản phẩm                                                                              *b Viết hàm tạo danh sách sản phẩm`n*c Viết hàm nhập các sản phẩm vào danh sách (từ File)`n*d Viết hàm xuất danh sách các sản phẩm (xuất ra đầy đủ Mã sản phẩm, Tên sản phẩm, Giá bán, Số lượng, Ngày sản xuất, Hạn sử dụng, Nhà cung cấp)`n*e Tạo Menu gồm các chức năng: Xuất danh sách sản phẩm, Thêm sản phẩm mới, Bổ sung số lượng sản phẩm, Xóa sản phẩm khỏi danh sách, Tìm kiếm sản phẩm, Bán sản phẩm, Sắp xếp sản phẩm theo giá bán, Tổng doanh thu bán hàng, Kiểm tra các sản phẩm đã được bán kèm theo sản phẩm bán chạy nhất, Kiểm tra hạn sử dụng các mặt hàng, Lợi nhuận`n*f Viết hàm hủy danh sách"
$ws.Range("B5").Value = $b5

# --- Update task description text (row 8 / col B): purchase/invoice description ---
$b8 = "*a Viết chức năng Mua sản phẩm (Yêu cầu nhập Tên sản phẩm)`n*b Viết chức năng Xuất hóa đơn (Xuất ra Tên các sản phẩm đã mua, số lượng mỗi sản phẩm, giá bán mỗi sản phẩm, thành tiền, tiền giảm, tổng giá đơn hàng)`n*c Viết chức năng tính Tổng doanh thu bán hàng     "
$ws.Range("B8").Value = $b8

# --- Fill in newly completed actual start/end dates (progress update) ---
$ws.Range("H11").Value = 43403
$ws.Range("G12").Value = 43408
$ws.Range("H12").Value = 43411
$ws.Range("G13").Value = 43413

# --- Update the view: reset scrolled top-left cell back to A1, move selection to H15 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H15").Select() | Out-Null
